$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.499.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.10%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.807.59"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.17%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.24%  "

$ws.Range("E6").Value = "  +4.23%  "

$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "38.87"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.19%  "

$ws.Range("E9").Value = "  -3.12%  "

$ws.Range("E10").Value = "  -3.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0986"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.30%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.067.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.09"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.23%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.808.50"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.03%  "

$ws.Range("E15").Value = "  -2.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.480.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.16%  "

$ws.Range("E17").Value = "  -2.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.68%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.85%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0770"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.73%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.46%  "

$ws.Range("E22").Value = "  -0.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.86%  "

$ws.Range("E24").Value = "  +1.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.53%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.61%  "

$ws.Range("E28").Value = "  +3.79%  "

$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.21%  "

$ws.Range("E31").Value = "  -2.17%  "

$ws.Range("E32").Value = "  -2.67%  "

$ws.Range("E33").Value = "  -4.15%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.94%  "

$ws.Range("E35").Value = "  -4.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.308.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.00%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.06"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.50%  "

$ws.Range("E38").Value = "  -1.81%  "

$ws.Range("E39").Value = "  -5.20%  "

$ws.Range("E40").Value = "  -0.25%  "

$ws.Range("E41").Value = "  +0.67%  "

$ws.Range("E42").Value = "  -0.49%  "

$ws.Range("E43").Value = "  +1.21%  "

$ws.Range("E44").Value = "  -1.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0514"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.35%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.968.73"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.89%  "

$ws.Range("E49").Value = "  -0.19%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.43%  "

$ws.Range("E51").Value = "  -6.38%  "
